# Update countries & provincias Spain
# Applies the periodic COVID data refresh:
#  - updates the "Datos actualizados" timestamp in A1
#  - updates rows whose country rankings/values changed after the data refresh
#    (Ucrania overtakes Singapur, Estado de Palestina overtakes Madagascar/Sri
#    Lanka/Guinea Ecuatorial, Fiyi/Dominica swap order, plus direct numeric
#    refreshes for Armenia, Hungria and Georgia)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 09:18"

function Set-Row($row, $pais, $b, $c, $d, $e, $f, $g, $h) {
    $ws.Cells.Item($row, 1).Value = $pais
    $ws.Cells.Item($row, 2).Value = $b
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $ws.Cells.Item($row, 5).Value = $e
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
}

# Ucrania now ranks above Singapur
Set-Row 37 "Ucrania"  43628 646 19027 23454 0 18 1147
Set-Row 38 "Singapur" 43459 0   37508 5925  0 0  26

# Armenia numbers refreshed (no rank change)
Set-Row 51 "Armenia" 25127 482 13297 11397 0 7 433

# Hungria numbers refreshed (no rank change)
Set-Row 93 "Hungria" 4145 3 2685 875 0 4 585

# Estado de Palestina now ranks above Madagascar, Sri Lanka, Guinea Ecuatorial
Set-Row 109 "Estado de Palestina" 2087 97 447  1635 0 1 5
Set-Row 110 "Madagascar"          2078 0  944  1116 0 0 18
Set-Row 111 "Sri Lanka"           2037 0  1661 365  0 0 11
Set-Row 112 "Guinea Ecuatorial"   2001 0  515  1454 0 0 32

# Georgia numbers refreshed (no rank change)
Set-Row 137 "Georgia" 926 2 791 120 0 0 15

# Fiyi / Dominica swap order (tied totals, values unchanged)
Set-Row 205 "Fiyi"     18 0 18 0 0 0 0
Set-Row 206 "Dominica" 18 0 18 0 0 0 0
